$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7143453.5
$ws.Range("I19").Value = 14286249
$ws.Range("K19").Value = 14286249
$ws.Range("M19").Value = -14286074

$ws.Range("H32").Value = 425.5

$ws.Range("H113").Value = 4253.3335
$ws.Range("I113").Value = 3944.375
$ws.Range("K113").Value = 3944.375
$ws.Range("M113").Value = -690.375

$ws.Range("H129").Value = 1067.48
$ws.Range("I129").Value = 363.22223
$ws.Range("J129").Value = 1137.1318
$ws.Range("K129").Value = 1089.66669
$ws.Range("L129").Value = 3411.3954
$ws.Range("M129").Value = 3910.33331
$ws.Range("N129").Value = -13411.3954

$ws.Range("H135").Value = 1073.9678
$ws.Range("I135").Value = 910.28
$ws.Range("J135").Value = 1756
$ws.Range("K135").Value = 8192.52
$ws.Range("L135").Value = 15804
$ws.Range("M135").Value = -5657.52
$ws.Range("N135").Value = -20874

$ws.Range("H137").Value = 1788018
$ws.Range("I137").Value = 2224398.2
$ws.Range("J137").Value = 2826.182
$ws.Range("K137").Value = 6673194.600000001
$ws.Range("L137").Value = 8478.545999999998
$ws.Range("M137").Value = -6670644.600000001
$ws.Range("N137").Value = -13578.546

$ws.Range("H141").Value = 207048.08
$ws.Range("I141").Value = 964.1111
$ws.Range("K141").Value = 2892.3333
$ws.Range("M141").Value = 2287.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1393.6531
$ws.Range("I61").Value = 618.7857
$ws.Range("J61").Value = 6042.857
$ws.Range("K61").Value = 618.7857
$ws.Range("L61").Value = 6042.857
$ws.Range("M61").Value = -406.7857
$ws.Range("N61").Value = -6466.857

$ws.Range("H63").Value = 2940.5
$ws.Range("I63").Value = 1854
$ws.Range("J63").Value = 6200
$ws.Range("K63").Value = 1854
$ws.Range("L63").Value = 6200
$ws.Range("M63").Value = -1168
$ws.Range("N63").Value = -7572

$ws.Range("H64").Value = 28791.666
$ws.Range("J64").Value = 28791.666
$ws.Range("L64").Value = 28791.666
$ws.Range("N64").Value = -29287.666

$ws.Range("H66").Value = 2940.5
$ws.Range("I66").Value = 1854
$ws.Range("J66").Value = 6200
$ws.Range("K66").Value = 9270
$ws.Range("L66").Value = 31000
$ws.Range("M66").Value = -5838
$ws.Range("N66").Value = -37864

$ws.Range("H67").Value = 28791.666
$ws.Range("J67").Value = 28791.666
$ws.Range("L67").Value = 28791.666
$ws.Range("N67").Value = -30507.666

$ws.Range("H74").Value = 611.7037
$ws.Range("I74").Value = 623.6923
$ws.Range("K74").Value = 623.6923
$ws.Range("M74").Value = 250.3077

$ws.Range("H76").Value = 29000
$ws.Range("J76").Value = 29000
$ws.Range("L76").Value = 29000
$ws.Range("N76").Value = -29676

$ws.Range("H77").Value = 611.7037
$ws.Range("I77").Value = 623.6923
$ws.Range("K77").Value = 3118.4615
$ws.Range("M77").Value = 1249.5385

$ws.Range("H79").Value = 29000
$ws.Range("J79").Value = 29000
$ws.Range("L79").Value = 29000
$ws.Range("N79").Value = -31340

$ws.Range("H119").Value = 739613.6
$ws.Range("J119").Value = 739613.6
$ws.Range("L119").Value = 739613.6
$ws.Range("N119").Value = -749289.6

$ws.Range("H132").Value = 1986.8937
$ws.Range("I132").Value = 1470.2632
$ws.Range("J132").Value = 4168.222
$ws.Range("K132").Value = 4410.7896
$ws.Range("L132").Value = 12504.666
$ws.Range("M132").Value = -1880.7896
$ws.Range("N132").Value = -17564.666

$ws.Range("H136").Value = 1393.6531
$ws.Range("I136").Value = 618.7857
$ws.Range("J136").Value = 6042.857
$ws.Range("K136").Value = 1856.3571
$ws.Range("L136").Value = 18128.571
$ws.Range("M136").Value = 693.6428999999998
$ws.Range("N136").Value = -23228.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 31000
$ws.Range("J63").Value = 31000
$ws.Range("L63").Value = 31000
$ws.Range("N63").Value = -32372

$ws.Range("H66").Value = 31000
$ws.Range("J66").Value = 31000
$ws.Range("L66").Value = 93000
$ws.Range("N66").Value = -99864

$ws.Range("H134").Value = 2151.5918
$ws.Range("I134").Value = 1595.4286
$ws.Range("J134").Value = 5488.5713
$ws.Range("K134").Value = 4786.2858
$ws.Range("L134").Value = 16465.7139
$ws.Range("M134").Value = -2251.2858
$ws.Range("N134").Value = -21535.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2312.5557
$ws.Range("I16").Value = 1725
$ws.Range("J16").Value = 2782.6
$ws.Range("K16").Value = 1725
$ws.Range("L16").Value = 2782.6
$ws.Range("M16").Value = -1438
$ws.Range("N16").Value = -3356.6

$ws.Range("H31").Value = 1925659.9
$ws.Range("I31").Value = 2942662
$ws.Range("J31").Value = 4655.6113
$ws.Range("K31").Value = 2942662
$ws.Range("L31").Value = 4655.6113
$ws.Range("M31").Value = -2942367
$ws.Range("N31").Value = -5245.6113

$ws.Range("H34").Value = 1925659.9
$ws.Range("I34").Value = 2942662
$ws.Range("J34").Value = 4655.6113
$ws.Range("K34").Value = 2942662
$ws.Range("L34").Value = 4655.6113
$ws.Range("M34").Value = -2942460
$ws.Range("N34").Value = -5059.6113

$ws.Range("H58").Value = 7144658.5
$ws.Range("I58").Value = 897.58185
$ws.Range("K58").Value = 897.58185
$ws.Range("M58").Value = -694.58185

$ws.Range("H113").Value = 2312.5557
$ws.Range("I113").Value = 1725
$ws.Range("J113").Value = 2782.6
$ws.Range("K113").Value = 1725
$ws.Range("L113").Value = 2782.6
$ws.Range("M113").Value = 445
$ws.Range("N113").Value = -7122.6

$ws.Range("H132").Value = 1974.3422
$ws.Range("I132").Value = 1401.4546
$ws.Range("K132").Value = 4204.3638
$ws.Range("M132").Value = -1674.3638

$ws.Range("H134").Value = 1273.76
$ws.Range("I134").Value = 680.12195
$ws.Range("K134").Value = 2040.36585
$ws.Range("M134").Value = 494.6341500000001

$ws.Range("H136").Value = 7144658.5
$ws.Range("I136").Value = 897.58185
$ws.Range("K136").Value = 2692.74555
$ws.Range("M136").Value = -142.7455500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 226.61539
$ws.Range("I12").Value = 31.375
$ws.Range("J12").Value = 313.3889
$ws.Range("K12").Value = 94.125
$ws.Range("L12").Value = 940.1667
$ws.Range("M12").Value = 78.875
$ws.Range("N12").Value = -1286.1667

$ws.Range("H68").Value = 3739.5
$ws.Range("I68").Value = 696
$ws.Range("J68").Value = 4500.375
$ws.Range("K68").Value = 2088
$ws.Range("L68").Value = 13501.125
$ws.Range("M68").Value = -1277
$ws.Range("N68").Value = -15123.125

$ws.Range("H71").Value = 3739.5
$ws.Range("I71").Value = 696
$ws.Range("J71").Value = 4500.375
$ws.Range("K71").Value = 6264
$ws.Range("L71").Value = 40503.375
$ws.Range("M71").Value = -2208
$ws.Range("N71").Value = -48615.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2561.561
$ws.Range("I132").Value = 2187.5
$ws.Range("J132").Value = 3581.7273
$ws.Range("K132").Value = 6562.5
$ws.Range("L132").Value = 10745.1819
$ws.Range("M132").Value = -4032.5
$ws.Range("N132").Value = -15805.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 200002500
$ws.Range("J22").Value = 4000
$ws.Range("L22").Value = 4000
$ws.Range("N22").Value = -4590

$ws.Range("H27").Value = 200002500
$ws.Range("J27").Value = 4000
$ws.Range("L27").Value = 4000
$ws.Range("N27").Value = -4214

$ws.Range("H111").Value = 39800
$ws.Range("J111").Value = 39800
$ws.Range("L111").Value = 39800
$ws.Range("N111").Value = -47980

$ws.Range("H132").Value = 1955.3256
$ws.Range("I132").Value = 1309.3
$ws.Range("J132").Value = 3446.1538
$ws.Range("K132").Value = 3927.9
$ws.Range("L132").Value = 10338.4614
$ws.Range("M132").Value = -1397.9
$ws.Range("N132").Value = -15398.4614

$ws.Range("H135").Value = 29314.5
$ws.Range("J135").Value = 29314.5
$ws.Range("L135").Value = 29314.5
$ws.Range("N135").Value = -39454.5

$ws.Range("H136").Value = 2327364
$ws.Range("I136").Value = 2858447.2
$ws.Range("K136").Value = 8575341.600000001
$ws.Range("M136").Value = -8572791.600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H119").Value = 29000
$ws.Range("J119").Value = 29000
$ws.Range("L119").Value = 29000
$ws.Range("N119").Value = -38676

$ws.Range("H132").Value = 153995.3
$ws.Range("I132").Value = 205944.33
$ws.Range("J132").Value = 26720.2
$ws.Range("K132").Value = 617832.99
$ws.Range("L132").Value = 80160.60000000001
$ws.Range("M132").Value = -615302.99
$ws.Range("N132").Value = -85220.60000000001

$ws.Range("H136").Value = 1113.1041
$ws.Range("I136").Value = 785.8
$ws.Range("K136").Value = 2357.4
$ws.Range("M136").Value = 192.6000000000004
